$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.16494083404541
$ws.Range("B1").Value = 2.548657655715942
$ws.Range("C1").Value = 1.170292019844055
$ws.Range("D1").Value = 0.5840051770210266
$ws.Range("E1").Value = 0.4197693467140198
